$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no numeric auto-conversion) for Price column cells
# whose new values are valid numeric literals, matching the workbook's
# original inline-string "Price" formatting.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D17", "D19", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = '27.983.46'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.857.31'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '312.37'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.5132'
$ws.Range("E7").Value = '  +1.26%  '
$ws.Range("D8").Value = '0.3831'
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '0.08229'
$ws.Range("E9").Value = '  -8.27%  '
$ws.Range("D10").Value = '1.108'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").Value = '41.52'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '6.182'
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("D13").Value = '20.54'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").Value = '1.861.86'
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").Value = '7.241'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("D19").Value = '0.06645'
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("D20").Value = '17.66'
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = '5.995'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '28.009.04'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").Value = '11.05'
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("D25").Value = '2.244'
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("D26").Value = '2.076.85'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = '2.507'
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").Value = '157.98'
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D29").Value = '20.42'
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("D30").Value = '124.37'
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").Value = '0.1063'
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("E32").Value = '  -2.78%  '
$ws.Range("D33").Value = '5.973'
$ws.Range("E33").Value = '  +6.53%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '9.365'
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("D36").Value = '0.02413'
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("D37").Value = '0.06484'
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("D38").Value = '0.2169'
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("D39").Value = '0.6517'
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("D40").Value = '1.193'
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D41").Value = '5.024'
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("D42").Value = '1.219'
$ws.Range("E42").Value = '  -3.64%  '
$ws.Range("D43").Value = '11.13'
$ws.Range("E43").Value = '  -2.53%  '
$ws.Range("D44").Value = '0.6137'
$ws.Range("E44").Value = '  +2.11%  '
$ws.Range("D45").Value = '13.00'
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").Value = '3.664'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '2.005'
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").Value = '120.22'
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").Value = '78.25'
$ws.Range("E51").Value = '  -1.73%  '
